# "Updates to final report"
# Build out Sheet2 as the condensed/final expense report (mirrors the
# Expense/Amount columns of Sheet1), then switch the active tab to it.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2 data -----------------------------------------------------
# Title
$ws2.Range("A1").Value = "Expense"
$ws2.Range("A1").Font.Size = 16
$ws2.Rows.Item(1).RowHeight = 21

# Header row
$ws2.Range("A2").Value = "Type"
$ws2.Range("B2").Value = "Amount"

# Expense rows (name, amount)
$ws2.Range("A3").Value = "Kinect"
$ws2.Range("B3").Value = 146.89

$ws2.Range("A4").Value = "Sonar Sensors"
$ws2.Range("B4").Value = 138.35

$ws2.Range("A5").Value = "PCBs"
$ws2.Range("B5").Value = 144.07

$ws2.Range("A6").Value = "Sheet Metal"
$ws2.Range("B6").Value = 50

$ws2.Range("A7").Value = "Kinect Mount"
$ws2.Range("B7").Value = 30

$ws2.Range("A8").Value = "Sonar Customs"
$ws2.Range("B8").Value = 20

$ws2.Range("A9").Value = "PCB Customs"
$ws2.Range("B9").Value = 21.65

$ws2.Range("A10").Value = "Paint"
$ws2.Range("B10").Value = 17.05

$ws2.Range("A11").Value = "Parts for PCB"
$ws2.Range("B11").Value = 85.76

$ws2.Range("A12").Value = "Display Board"
$ws2.Range("B12").Value = 67

# Total row
$ws2.Range("A15").Value = "Total"
$ws2.Range("B15").Formula = "=SUM(B3:B12)"
$ws2.Range("A15:B15").Font.Bold = $true

# Column widths (best-fit sized, as the original report's columns were)
$ws2.Columns.Item(1).ColumnWidth = 13.3
$ws2.Columns.Item(2).ColumnWidth = 7.3

# --- View / selection changes -----------------------------------------
# Sheet1 selection moves from the single cell E4 to the totals block,
# and Sheet1 is no longer the tab shown when the workbook opens.
$ws1.Range("A1:D15").Select() | Out-Null

# Sheet2 becomes the active/visible tab.
$ws2.Activate() | Out-Null
$ws2.Range("A1").Select() | Out-Null
